$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '60.945.60'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  +7.01%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.674.29'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  +9.99%  '
$ws.Range('E4').Value = '  -0.10%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '511.86'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +4.89%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '157.52'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +2.67%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.998'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -0.20%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.605'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +0.14%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '2.669.88'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +9.98%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '6.37'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +7.31%  '
$ws.Range('E11').Value = '  +5.05%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.349'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +4.04%  '
$ws.Range('E13').Value = '  +1.11%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '3.124.12'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +9.67%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '60.946.83'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +6.77%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '21.81'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +4.75%  '
$ws.Range('E17').Value = '  +4.74%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '2.667.71'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +9.77%  '
$ws.Range('E19').Value = '  +0.86%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '348.39'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +7.13%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '10.50'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +4.91%  '
$ws.Range('E22').Value = '  +3.56%  '
$ws.Range('E23').Value = '  -0.19%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '60.23'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +3.53%  '
$ws.Range('E25').Value = '  +3.12%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.778.57'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +9.71%  '
$ws.Range('E27').Value = '  +3.41%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.995'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -0.44%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '0.0₃0862'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +9.81%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '7.53'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +2.75%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.999'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -0.10%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '157.65'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +5.15%  '
$ws.Range('E33').Value = '  +4.96%  '
$ws.Range('E34').Value = '  +3.43%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '5.70'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +6.29%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '4.05'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +8.72%  '
$ws.Range('E37').Value = '  +5.14%  '
$ws.Range('E38').Value = '  +11.09%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '309.60'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +15.25%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.861'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +1.39%  '
$ws.Range('E41').Value = '  +6.30%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.837'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +28.56%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '35.43'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +3.89%  '
$ws.Range('E44').Value = '  +8.39%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.0578'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +8.33%  '
$ws.Range('E46').Value = '  -0.68%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.999'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -0.01%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '19.88'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +14.16%  '
$ws.Range('B49').Value = 'RenderToken'
$ws.Range('C49').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '4.86'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +5.55%  '
$ws.Range('B50').Value = 'VeChain'
$ws.Range('C50').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.0236'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +3.50%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '2.043.90'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +9.30%  '
